# Generate Report for Handback
# Refresh the handoff/handback timestamps for the "4bed097a-..." report row
# (the zh-cn and de-de sheets), and roll the latest of those two new
# timestamps up into the Overview sheet's "Latest HO Xliff Generate Date"
# column.

$wb = $excel.ActiveWorkbook

$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("H2").Value = "2016-08-30 08:43:45"
$zh.Range("K2").Value = "2016-08-30 08:44:40"

$de = $wb.Worksheets.Item("de-de")
$de.Range("H2").Value = "2016-08-30 08:43:57"
$de.Range("K2").Value = "2016-08-30 08:44:57"

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G2").Value = "2016-08-30 08:43:57"
